$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Rows 74 / 75 (match index A74=73, A75=74): the two fixtures'
#    details (F:V) were swapped - Mafra-Leixoes now sits on row 75
#    and FC Porto B-Feirense now sits on row 74.
# -----------------------------------------------------------------
$ws.Range("F74").Value = "FC Porto B"
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = "Feirense"
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1.98
$ws.Range("K74").Value = "04/11/2023 11:48"
$ws.Range("L74").Value = 1.84
$ws.Range("M74").Value = "04/11/2023 11:59"
$ws.Range("N74").Value = 3.57
$ws.Range("O74").Value = "04/11/2023 11:48"
$ws.Range("P74").Value = 3.78
$ws.Range("Q74").Value = "04/11/2023 11:59"
$ws.Range("R74").Value = 3.87
$ws.Range("S74").Value = "04/11/2023 11:48"
$ws.Range("T74").Value = 4.41
$ws.Range("U74").Value = "04/11/2023 11:58"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-feirense/jTL6QSDN/"

$ws.Range("F75").Value = "Mafra"
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = "Leixoes"
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = 1.88
$ws.Range("K75").Value = "04/11/2023 11:48"
$ws.Range("L75").Value = 1.93
$ws.Range("M75").Value = "04/11/2023 11:48"
$ws.Range("N75").Value = 3.73
$ws.Range("O75").Value = "04/11/2023 11:48"
$ws.Range("P75").Value = 3.55
$ws.Range("Q75").Value = "04/11/2023 11:51"
$ws.Range("R75").Value = 3.8
$ws.Range("S75").Value = "04/11/2023 11:48"
$ws.Range("T75").Value = 4.21
$ws.Range("U75").Value = "04/11/2023 11:51"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leixoes/YiBBPnTT/"

# -----------------------------------------------------------------
# 2) Rows 88 / 89 (match index A88=87, A89=88): likewise swapped -
#    Torreense-Mafra now sits on row 88 and Vilaverdense-Tondela now
#    sits on row 89.
# -----------------------------------------------------------------
$ws.Range("F88").Value = "Torreense"
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = "Mafra"
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2.12
$ws.Range("K88").Value = "08/11/2023 06:12"
$ws.Range("L88").Value = 2.57
$ws.Range("M88").Value = "12/11/2023 11:50"
$ws.Range("N88").Value = 3.41
$ws.Range("O88").Value = "08/11/2023 06:12"
$ws.Range("P88").Value = 3.18
$ws.Range("Q88").Value = "12/11/2023 11:50"
$ws.Range("R88").Value = 3.58
$ws.Range("S88").Value = "08/11/2023 06:12"
$ws.Range("T88").Value = 3.04
$ws.Range("U88").Value = "12/11/2023 11:50"
$ws.Range("V88").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/torreense-mafra/jTsxt8DA/"

$ws.Range("F89").Value = "Vilaverdense"
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = "Tondela"
$ws.Range("I89").Value = 2
$ws.Range("J89").Value = 2.68
$ws.Range("K89").Value = "08/11/2023 06:12"
$ws.Range("L89").Value = 3.73
$ws.Range("M89").Value = "12/11/2023 11:57"
$ws.Range("N89").Value = 3.34
$ws.Range("O89").Value = "08/11/2023 06:12"
$ws.Range("P89").Value = 3.49
$ws.Range("Q89").Value = "12/11/2023 11:57"
$ws.Range("R89").Value = 2.58
$ws.Range("S89").Value = "08/11/2023 06:12"
$ws.Range("T89").Value = 2.09
$ws.Range("U89").Value = "12/11/2023 11:57"
$ws.Range("V89").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/vilaverdense-fc-tondela/UJttuSSG/"

# -----------------------------------------------------------------
# 3) New row 104 appended at the bottom (Nacional 2-1 Pacos Ferreira).
#    Copy the formatting (bold/border/center for A, date number format
#    for E) from the previous last row (103) before writing the values
#    so the new row matches the sheet's existing look.
# -----------------------------------------------------------------
$ws.Range("A103").Copy() | Out-Null
$ws.Range("A104").PasteSpecial(-4122) | Out-Null

$ws.Range("E103").Copy() | Out-Null
$ws.Range("E104").PasteSpecial(-4122) | Out-Null

$ws.Range("A104").Value = 103
$ws.Range("B104").Value = "portugal"
$ws.Range("C104").Value = "liga-portugal-2"
$ws.Range("D104").Value = "2023-2024"
$ws.Range("E104").Value = 45262.625
$ws.Range("F104").Value = "Nacional"
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = "Pacos Ferreira"
$ws.Range("I104").Value = 1
$ws.Range("J104").Value = 2.12
$ws.Range("K104").Value = "27/11/2023 13:12"
$ws.Range("L104").Value = 1.99
$ws.Range("M104").Value = "02/12/2023 14:58"
$ws.Range("N104").Value = 3.45
$ws.Range("O104").Value = "27/11/2023 13:12"
$ws.Range("P104").Value = 3.69
$ws.Range("Q104").Value = "02/12/2023 14:58"
$ws.Range("R104").Value = 3.54
$ws.Range("S104").Value = "27/11/2023 13:12"
$ws.Range("T104").Value = 3.81
$ws.Range("U104").Value = "02/12/2023 14:58"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-pacos-ferreira/tz8ffbYK/"
